$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.689.92'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '1.813.83'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '''228.40'
$ws.Range("E5").Value = '  +0.88%  '
$ws.Range("D6").Value = '''0.566'
$ws.Range("E6").Value = '  +1.88%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").Value = '''34.84'
$ws.Range("E8").Value = '  +7.37%  '
$ws.Range("D9").Value = '''0.299'
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("D10").Value = '''0.0695'
$ws.Range("E10").Value = '  +0.49%  '
$ws.Range("D11").Value = '''0.0952'
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D12").Value = '2.077.42'
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").Value = '''11.30'
$ws.Range("E13").Value = '  +2.30%  '
$ws.Range("D14").Value = '1.823.26'
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("D15").Value = '''0.643'
$ws.Range("E15").Value = '  +2.03%  '
$ws.Range("D16").Value = '34.667.36'
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").Value = '''4.33'
$ws.Range("E17").Value = '  +2.56%  '
$ws.Range("D18").Value = '''69.10'
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("D19").Value = '''247.73'
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").Value = '0.0₃0800'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = '''11.54'
$ws.Range("E21").Value = '  +5.15%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '''4.18'
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("D24").Value = '''172.09'
$ws.Range("E24").Value = '  +6.13%  '
$ws.Range("E25").Value = '  +1.99%  '
$ws.Range("D26").Value = '''7.44'
$ws.Range("E26").Value = '  +3.49%  '
$ws.Range("D27").Value = '''16.75'
$ws.Range("E27").Value = '  +2.17%  '
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '''0.0532'
$ws.Range("E30").Value = '  +2.09%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''3.98'
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("D32").Value = '''3.85'
$ws.Range("E32").Value = '  +1.66%  '
$ws.Range("D34").Value = '''1.85'
$ws.Range("E34").Value = '  +1.63%  '
$ws.Range("D35").Value = '''2.66'
$ws.Range("E35").Value = '  +2.38%  '
$ws.Range("D36").Value = '1.419.20'
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("D37").Value = '''0.679'
$ws.Range("E37").Value = '  +2.25%  '
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("D39").Value = '''0.0191'
$ws.Range("E39").Value = '  +0.54%  '
$ws.Range("D40").Value = '''85.64'
$ws.Range("E40").Value = '  +2.81%  '
$ws.Range("D41").Value = '''2.84'
$ws.Range("E41").Value = '  +3.61%  '
$ws.Range("D42").Value = '''0.957'
$ws.Range("E42").Value = '  +2.80%  '
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("D44").Value = '''13.87'
$ws.Range("E44").Value = '  -1.16%  '
$ws.Range("D45").Value = '''0.0523'
$ws.Range("E45").Value = '  +0.43%  '
$ws.Range("E46").Value = '  +3.26%  '
$ws.Range("D47").Value = '''6.13'
$ws.Range("E47").Value = '  +0.83%  '
$ws.Range("D48").Value = '1.977.29'
$ws.Range("E48").Value = '  +1.52%  '
$ws.Range("D49").Value = '''105.58'
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("E51").Value = '  +0.22%  '
